# TC20_Canine_Filter_Breed-Dalmatian.xlsx — "Fixed ICDC breed all testcases"
#
# The "StatQuery" column (C) on the "startup" tab previously pointed every
# tab (Cases/Samples/Files) at a single combined counts query. That query is
# replaced with a corrected Cypher query that reports Programs/Studies/
# Cases/Samples/Case Files/Study Files, and it is now used for all three
# rows (CasesTab, SamplesTab, FilesTab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Dalmatian']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Rows 2-4 are CasesTab, SamplesTab and FilesTab respectively; column C
# ("StatQuery") gets the corrected query on every row.
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# The sheet view moved down one row and zoomed in slightly.
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B3").Select()
